# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.323.10'
$ws.Range("D3").Value = '1.871.99'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''243.64'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.4696'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").Value = '''0.2876'
$ws.Range("E8").Value = '  -0.87%  '
$ws.Range("D9").Value = '''0.06445'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").Value = '''22.11'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").Value = '''0.07769'
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").Value = '1.873.44'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D14").Value = '''0.7205'
$ws.Range("D15").Value = '''5.127'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '''279.48'
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("D17").Value = '30.315.24'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '''12.99'
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '''0.000007440'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '2.119.08'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '''5.232'
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").Value = '''6.233'
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("D25").Value = '''163.30'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").Value = '''9.050'
$ws.Range("D27").Value = '''18.67'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").Value = '''1.317'
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("D30").Value = '''0.09572'
$ws.Range("E30").Value = '  -3.22%  '
$ws.Range("D31").Value = '''1.469'
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("D32").Value = '''4.215'
$ws.Range("E32").Value = '  -0.69%  '
$ws.Range("D33").Value = '''4.092'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").Value = '''0.04805'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("D36").Value = '''0.6874'
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").Value = '''0.01870'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '''2.813'
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").Value = '''6.228'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '''74.25'
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").Value = '''0.4235'
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("D43").Value = '''1.935'
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = '''0.9993'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '''0.8241'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Value = '''100.75'
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("D47").Value = '''9.554'
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("D48").Value = '''35.11'
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").Value = '''6.916'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").Value = '''899.44'
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("D51").Value = '''0.05717'
$ws.Range("E51").Value = '  +0.89%  '
